$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Sending cluster) changes from "Resolving-Mac" to "Inflammatory-Mac" for all data rows
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("A6").Value = "Inflammatory-Mac"

# Updated TPM-derived numeric values for rows 2-6 (columns E:T)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7878926666666667
$ws.Range("H2").Value = 2.363678
$ws.Range("M2").Value = 12.86269466666666
$ws.Range("N2").Value = 38.58808399999999
$ws.Range("O2").Value = 0.1337831063410017
$ws.Range("P2").Value = 0.1337831063410017
$ws.Range("Q2").Value = 10.13442280143911
$ws.Range("R2").Value = 91.20980521295199
$ws.Range("S2").Value = 0.1337831063410017
$ws.Range("T2").Value = 0.1337831063410017

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7878926666666667
$ws.Range("H3").Value = 2.363678
$ws.Range("O3").Value = 0.3593152390330854
$ws.Range("P3").Value = 0.3593152390330854
$ws.Range("Q3").Value = 27.21907609231089
$ws.Range("R3").Value = 244.971684830798
$ws.Range("S3").Value = 0.3593152390330854
$ws.Range("T3").Value = 0.3593152390330854

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7878926666666667
$ws.Range("H4").Value = 2.363678
$ws.Range("M4").Value = 18.65324433333334
$ws.Range("N4").Value = 55.95973300000001
$ws.Range("O4").Value = 0.1940098117012772
$ws.Range("P4").Value = 0.1940098117012772
$ws.Range("Q4").Value = 14.69675441977489
$ws.Range("R4").Value = 132.270789777974
$ws.Range("S4").Value = 0.1940098117012772
$ws.Range("T4").Value = 0.1940098117012772

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7878926666666667
$ws.Range("H5").Value = 2.363678
$ws.Range("M5").Value = 7.643308666666666
$ws.Range("N5").Value = 22.929926
$ws.Range("O5").Value = 0.07949699519803316
$ws.Range("P5").Value = 0.07949699519803316
$ws.Range("Q5").Value = 6.022106847536445
$ws.Range("R5").Value = 54.198961627828
$ws.Range("S5").Value = 0.07949699519803316
$ws.Range("T5").Value = 0.07949699519803316

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7878926666666667
$ws.Range("H6").Value = 2.363678
$ws.Range("M6").Value = 22.43995333333334
$ws.Range("N6").Value = 67.31986000000001
$ws.Range("O6").Value = 0.2333948477266026
$ws.Range("P6").Value = 0.2333948477266026
$ws.Range("Q6").Value = 17.68027467167556
$ws.Range("R6").Value = 159.12247204508
$ws.Range("S6").Value = 0.2333948477266026
$ws.Range("T6").Value = 0.2333948477266026
